$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update value in B8 with new result
$ws.Range("B8").Value = 0.3417

# Move active cell/view selection to C1 (no more topLeftCell scroll to A7)
$ws.Range("C1").Select()
